# Regenerate column G ("K") values for rows 2-14 on Sheet1,
# per the diff: column K (G) recalculated to use "K" instead of "Strike#".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 3
    3  = 0
    4  = 1
    5  = 0
    6  = 3
    7  = 2
    8  = 0
    9  = 1
    10 = 0
    11 = 4
    12 = 0
    13 = 0
    14 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
